# Update countries & provincias Spain
# - Swap the sorted-order positions of "Grecia" / "Republica de Yibuti"
#   (rows 102/103) and "Santa Lucia" / "Timor Oriental" (rows 202/203) -
#   the underlying country names for those two row-pairs trade places.
# - Refresh the "datos actualizados" timestamp string.
# - Update the numeric statistics (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Casos criticos, Muertes) for the rows whose
#   counts changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Country-name swaps ------------------------------------------------
$ws.Range("A102").Value = "Grecia"
$ws.Range("A103").Value = "Republica de Yibuti"

$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Timestamp string ----------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 18:13"

# --- Row 4: Estados Unidos ------------------------------------------------
$ws.Range("B4").Value = 5111473
$ws.Range("C4").Value = 15949
$ws.Range("D4").Value = 2618440
$ws.Range("E4").Value = 2328572
$ws.Range("G4").Value = 367
$ws.Range("H4").Value = 164461

# --- Row 6: India -----------------------------------------------------
$ws.Range("B6").Value = 2129154
$ws.Range("C6").Value = 42290
$ws.Range("D6").Value = 1461772
$ws.Range("E6").Value = 624238
$ws.Range("G6").Value = 566
$ws.Range("H6").Value = 43144

# --- Row 15: Reino Unido ------------------------------------------------
$ws.Range("G15").Value = 55
$ws.Range("H15").Value = 46566

# --- Row 19: Italia -----------------------------------------------------
$ws.Range("B19").Value = 250103
$ws.Range("C19").Value = 347
$ws.Range("D19").Value = 201947
$ws.Range("E19").Value = 12953
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 35203

# --- Row 24: Irak -------------------------------------------------------
$ws.Range("B24").Value = 147389
$ws.Range("C24").Value = 3325
$ws.Range("D24").Value = 105504
$ws.Range("E24").Value = 36575
$ws.Range("G24").Value = 74
$ws.Range("H24").Value = 5310

# --- Row 74: Chequia ------------------------------------------------------
$ws.Range("B74").Value = 18146
$ws.Range("C74").Value = 86
$ws.Range("D74").Value = 12764
$ws.Range("E74").Value = 4993

# --- Row 99: Albania ------------------------------------------------------
$ws.Range("B99").Value = 6275
$ws.Range("C99").Value = 124
$ws.Range("D99").Value = 3268
$ws.Range("E99").Value = 2814
$ws.Range("G99").Value = 4
$ws.Range("H99").Value = 193

# --- Row 102: now Grecia (updated data) ------------------------------
$ws.Range("B102").Value = 5421
$ws.Range("C102").Value = 151
$ws.Range("D102").Value = 1374
$ws.Range("E102").Value = 3836
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 211

# --- Row 103: now Republica de Yibuti (carries old Grecia-row totals) --
$ws.Range("B103").Value = 5338
$ws.Range("D103").Value = 5083
$ws.Range("E103").Value = 196
$ws.Range("H103").Value = 59

# --- Row 119: Sri Lanka ---------------------------------------------------
$ws.Range("B119").Value = 2841
$ws.Range("C119").Value = 2
$ws.Range("E119").Value = 254

# --- Row 143: Jordania ----------------------------------------------------
$ws.Range("B143").Value = 1246
$ws.Range("C143").Value = 9
$ws.Range("E143").Value = 57

# --- Row 179: Trinidad yTobago --------------------------------------------
$ws.Range("B179").Value = 243
$ws.Range("C179").Value = 18
$ws.Range("E179").Value = 100
